$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells we touch to remain Text (some look like numbers/dates,
# e.g. "27.101.49" or "0.09210" which would otherwise lose precision/format
# if Excel auto-converted them to a numeric value).
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated price (D) and volume-change (E) values scraped for this run.
$ws.Range("D2").Value = "27.101.49"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.900.62"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "307.22"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "0.5235"
$ws.Range("E7").Value = "  +3.15%  "
$ws.Range("D8").Value = "0.3778"
$ws.Range("E8").Value = "  +3.06%  "
$ws.Range("D9").Value = "0.07225"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "21.16"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "0.8921"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "0.07678"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").Value = "1.907.87"
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").Value = "94.36"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").Value = "5.233"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "0.000008517"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "14.51"
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("D19").Value = "0.9997"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "27.140.07"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "5.067"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").Value = "2.139.95"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("D23").Value = "10.61"
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("D24").Value = "6.414"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "2.294"
$ws.Range("E25").Value = "  +9.67%  "
$ws.Range("D26").Value = "145.75"
$ws.Range("E26").Value = "  -1.88%  "
$ws.Range("D27").Value = "1.735"
$ws.Range("E27").Value = "  -2.96%  "
$ws.Range("D28").Value = "18.09"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").Value = "114.63"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").Value = "4.963"
$ws.Range("E30").Value = "  +4.49%  "
$ws.Range("D31").Value = "4.797"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").Value = "0.09210"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").Value = "0.05054"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").Value = "1.239"
$ws.Range("E34").Value = "  +6.90%  "
$ws.Range("D35").Value = "0.7771"
$ws.Range("E35").Value = "  +3.46%  "
$ws.Range("D36").Value = "2.976"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "3.298"
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").Value = "2.593"
$ws.Range("E38").Value = "  +1.91%  "
$ws.Range("D39").Value = "0.5649"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").Value = "0.01992"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").Value = "1.072"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").Value = "8.982"
$ws.Range("E42").Value = "  +4.75%  "
$ws.Range("D43").Value = "6.629"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "118.68"
$ws.Range("E44").Value = "  +2.63%  "
$ws.Range("D45").Value = "0.1518"
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("D46").Value = "0.4836"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("D47").Value = "10.17"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "1.598"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("D50").Value = "37.49"
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").Value = "64.09"
$ws.Range("E51").Value = "  +1.61%  "
